$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "69.289.89"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +2.23%  "

# Row 3
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "3.381.57"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +1.40%  "

# Row 4
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  +0.03%  "

# Row 5
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "586.80"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.92%  "

# Row 6
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "179.47"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +1.72%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.06%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +1.00%  "

# Row 9
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.195"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +6.30%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  +1.28%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  +2.70%  "

# Row 12
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "0.0000282"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +3.17%  "

# Row 13
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "681.33"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -1.91%  "

# Row 14
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "8.61"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +2.00%  "

# Row 15
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "3.917.03"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +1.13%  "

# Row 16
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "69.253.77"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +2.19%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  +1.64%  "

# Row 18
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "3.387.86"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +1.64%  "

# Row 19
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "17.65"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +0.60%  "

# Row 20
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "11.25"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +1.68%  "

# Row 21
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "0.903"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +1.03%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "  -0.93%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  +0.77%  "

# Row 24
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "103.36"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +2.99%  "

# Row 25
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "3.92"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +0.25%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  +0.75%  "

# Row 27
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "9.60"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +0.60%  "

# Row 28
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "34.01"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +2.87%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  +1.53%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  -1.81%  "

# Row 31
$ws.Cells.Item(31, 2).Value = "Cosmos"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "11.12"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +0.96%  "

# Row 32
$ws.Cells.Item(32, 2).Value = "Bittensor"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "556.07"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -1.86%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  +0.69%  "

# Row 34
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "3.55"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +4.95%  "

# Row 35
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "58.64"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +2.23%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  -0.15%  "

# Row 37
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "3.673.68"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -0.97%  "

# Row 38
$ws.Cells.Item(38, 2).Value = "Kaspa"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "0.139"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +3.89%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "35.44"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +2.09%  "

# Row 40
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "3.28"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +3.59%  "

# Row 41
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "2.67"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +1.13%  "

# Row 42
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "0.0₃0698"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +3.37%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  +0.82%  "

# Row 44
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "0.0422"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +3.22%  "

# Row 45
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "3.31"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +0.61%  "

# Row 46
$ws.Cells.Item(46, 5).Value = "  -0.15%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  +1.00%  "

# Row 48
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "1.42"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +5.43%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  +0.02%  "

# Row 50
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "132.52"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +1.01%  "

# Row 51
$ws.Cells.Item(51, 5).Value = "  +3.52%  "
